$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(4, 6).Value = 5125
$ws.Cells.Item(5, 6).Value = 5125
$ws.Cells.Item(7, 6).Value = 156
$ws.Cells.Item(8, 6).Value = 208
$ws.Cells.Item(12, 6).Value = 8509
$ws.Cells.Item(13, 6).Value = 8509
$ws.Cells.Item(16, 6).Value = 620
$ws.Cells.Item(17, 6).Value = 2558
$ws.Cells.Item(19, 6).Value = 2317
$ws.Cells.Item(20, 6).Value = 9
$ws.Cells.Item(22, 6).Value = 2533
$ws.Cells.Item(23, 6).Value = 23
$ws.Cells.Item(25, 6).Value = 6455
$ws.Cells.Item(28, 6).Value = 140
$ws.Cells.Item(31, 6).Value = 6964
$ws.Cells.Item(32, 6).Value = 7
$ws.Cells.Item(34, 6).Value = 234
$ws.Cells.Item(37, 6).Value = 108
$ws.Cells.Item(38, 6).Value = 9
$ws.Cells.Item(42, 6).Value = 52
$ws.Cells.Item(48, 6).Value = 532
$ws.Cells.Item(49, 6).Value = 2536
$ws.Cells.Item(50, 6).Value = 85

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(3, 6).Value = 16
$ws.Cells.Item(4, 6).Value = 185
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(6, 6).Value = 76

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(3, 6).Value = 5125
$ws.Cells.Item(4, 6).Value = 5125
$ws.Cells.Item(6, 6).Value = 156
$ws.Cells.Item(7, 6).Value = 208
$ws.Cells.Item(11, 6).Value = 8509
$ws.Cells.Item(12, 6).Value = 8509
$ws.Cells.Item(15, 6).Value = 620
$ws.Cells.Item(16, 6).Value = 2558
$ws.Cells.Item(17, 6).Value = 185
$ws.Cells.Item(19, 6).Value = 2317
$ws.Cells.Item(20, 6).Value = 76
$ws.Cells.Item(21, 6).Value = 9
$ws.Cells.Item(22, 6).Value = 2533
$ws.Cells.Item(23, 6).Value = 23
$ws.Cells.Item(27, 6).Value = 6455
$ws.Cells.Item(30, 6).Value = 140
$ws.Cells.Item(33, 6).Value = 6964
$ws.Cells.Item(34, 6).Value = 7
$ws.Cells.Item(36, 6).Value = 234
$ws.Cells.Item(38, 6).Value = 108
$ws.Cells.Item(41, 6).Value = 52
$ws.Cells.Item(47, 6).Value = 532
$ws.Cells.Item(49, 6).Value = 2536
$ws.Cells.Item(50, 6).Value = 85
